# Refactor plan rules into modular components and add response logic
# Update the AIP All-Eligible projection results with the refreshed
# simulation output (participation / deferral / contribution projections).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2 = @{
        C = 9458
        D = 7574
        E = 0.8008035525481074
        F = 0.7991137370753324
        G = 0.1028346976498548
        H = 0.08217661953998734
        I = 39239770.164845
        J = 13420220.8417585
        L = 13420220.8417585
        M = 52659991.0066035
        N = 800122294.0972
        O = 782422487.0932001
        P = 0.0167727120476013
        Q = 0.01715214102756216
    }
    3 = @{
        C = 9637
        D = 7551
        E = 0.7835425962436443
        F = 0.7811918063314711
        G = 0.1028830618461131
        H = 0.08037140492447756
        I = 40538728.96122567
        J = 13799638.26741243
        L = 13799638.26741243
        M = 54338367.2286381
        N = 835784763.7231281
        O = 818304587.6990581
        P = 0.01651099525425647
        Q = 0.01686369412423168
    }
    4 = @{
        D = 7543
        E = 0.7665650406504065
        F = 0.7651653479407587
        G = 0.1029497547394936
        H = 0.07877358490566039
        I = 42014448.3095379
        J = 14217426.49631654
        L = 14217426.49631654
        M = 56231874.80585443
        N = 874054288.5903099
        O = 856605340.5843561
        P = 0.01626606800276291
        Q = 0.01659740585625785
    }
    5 = @{
        C = 10034
        D = 7534
        E = 0.7508471197927048
        F = 0.7493534911478019
        G = 0.1030222989116007
        H = 0.07720011935548041
        I = 43535019.83059579
        J = 14644484.75463068
        L = 14644484.75463068
        M = 58179504.58522647
        N = 913242019.3379748
        O = 895755913.8750015
        P = 0.01603571062712021
        Q = 0.01634874470577512
    }
    6 = @{
        C = 10228
        D = 7523
        E = 0.7355299178725069
        F = 0.7336649112541447
        G = 0.1030931809118703
        H = 0.07563584942461479
        I = 45146428.42506469
        J = 15078985.98419153
        L = 15078985.98419153
        M = 60225414.4092562
        N = 954929691.5795953
        O = 937337865.7064139
        P = 0.01579067665101987
        Q = 0.01608703386033319
    }
}

foreach ($row in $updates.Keys) {
    $cols = $updates[$row]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$row").Value = $cols[$col]
    }
}
